$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of Fruta / hortaliza dataset: each row (2-31, excluding the
# already-correct row 23) gets new Fecha/Volumen/Precio/Origen values as the
# underlying data was reshuffled across the week.

# Row 2
$ws.Range("D2").Value = 44466
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 13500
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 13750
$ws.Range("S2").Value = 6875

# Row 3
$ws.Range("D3").Value = 44445
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 7250

# Row 4
$ws.Range("D4").Value = 44658
$ws.Range("N4").Value = 6500
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6750
$ws.Range("S4").Value = 3375

# Row 5
$ws.Range("D5").Value = 44497
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 11500
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 11750
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 5875

# Row 6
$ws.Range("D6").Value = 44452
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 13500
$ws.Range("R6").Value = 'Provincia de Limarí'
$ws.Range("S6").Value = 6750

# Row 7
$ws.Range("D7").Value = 44446
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("S7").Value = 7250

# Row 8
$ws.Range("D8").Value = 44435
$ws.Range("M8").Value = 400
$ws.Range("N8").Value = 19500
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19750
$ws.Range("S8").Value = 9875

# Row 9
$ws.Range("D9").Value = 44651
$ws.Range("M9").Value = 400
$ws.Range("N9").Value = 6000
$ws.Range("O9").Value = 6500
$ws.Range("P9").Value = 6250
$ws.Range("R9").Value = 'Provincia de Linares'
$ws.Range("S9").Value = 3125

# Row 10
$ws.Range("D10").Value = 44659
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 6500
$ws.Range("O10").Value = 7000
$ws.Range("P10").Value = 6750
$ws.Range("R10").Value = 'Provincia de Linares'
$ws.Range("S10").Value = 3375

# Row 11
$ws.Range("D11").Value = 44468
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 13000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 13500
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 6750

# Row 12
$ws.Range("D12").Value = 44637
$ws.Range("M12").Value = 160
$ws.Range("R12").Value = 'Región de O''Higgins'

# Row 13
$ws.Range("D13").Value = 44463
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 13000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 13500
$ws.Range("S13").Value = 6750

# Row 14
$ws.Range("D14").Value = 44490
$ws.Range("N14").Value = 11500
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 11750
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 5875

# Row 15
$ws.Range("D15").Value = 44495
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 11000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 11500
$ws.Range("S15").Value = 5750

# Row 16
$ws.Range("D16").Value = 44498
$ws.Range("M16").Value = 240
$ws.Range("N16").Value = 11000
$ws.Range("O16").Value = 11500
$ws.Range("P16").Value = 11250
$ws.Range("S16").Value = 5625

# Row 17
$ws.Range("D17").Value = 44459
$ws.Range("N17").Value = 13000
$ws.Range("O17").Value = 14000
$ws.Range("P17").Value = 13500
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 6750

# Row 18
$ws.Range("D18").Value = 44462
$ws.Range("M18").Value = 140

# Row 19
$ws.Range("D19").Value = 44494
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 11500
$ws.Range("O19").Value = 12000
$ws.Range("P19").Value = 11750
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 5875

# Row 20
$ws.Range("D20").Value = 44448
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 14500
$ws.Range("S20").Value = 7250

# Row 21
$ws.Range("D21").Value = 44645
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 6000
$ws.Range("O21").Value = 6500
$ws.Range("P21").Value = 6250
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 3125

# Row 22
$ws.Range("D22").Value = 44630
$ws.Range("M22").Value = 240
$ws.Range("N22").Value = 6000
$ws.Range("O22").Value = 6500
$ws.Range("P22").Value = 6250
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 3125

# Row 24
$ws.Range("D24").Value = 44638
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 6000
$ws.Range("O24").Value = 6500
$ws.Range("P24").Value = 6250
$ws.Range("R24").Value = 'Provincia de Linares'
$ws.Range("S24").Value = 3125

# Row 25
$ws.Range("D25").Value = 44634
$ws.Range("M25").Value = 160
$ws.Range("N25").Value = 6000
$ws.Range("O25").Value = 6500
$ws.Range("P25").Value = 6250
$ws.Range("R25").Value = 'Provincia de Linares'
$ws.Range("S25").Value = 3125

# Row 26
$ws.Range("D26").Value = 44665
$ws.Range("N26").Value = 6500
$ws.Range("O26").Value = 7000
$ws.Range("P26").Value = 6750
$ws.Range("R26").Value = 'Provincia de Linares'
$ws.Range("S26").Value = 3375

# Row 27
$ws.Range("D27").Value = 44644
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = 6000
$ws.Range("O27").Value = 6500
$ws.Range("P27").Value = 6250
$ws.Range("R27").Value = 'Región de O''Higgins'
$ws.Range("S27").Value = 3125

# Row 28
$ws.Range("D28").Value = 44491
$ws.Range("M28").Value = 200

# Row 29
$ws.Range("D29").Value = 44489
$ws.Range("M29").Value = 400
$ws.Range("N29").Value = 11500
$ws.Range("O29").Value = 12000
$ws.Range("P29").Value = 11750
$ws.Range("R29").Value = 'Provincia de Limarí'
$ws.Range("S29").Value = 5875

# Row 30
$ws.Range("D30").Value = 44631
$ws.Range("M30").Value = 160
$ws.Range("N30").Value = 6000
$ws.Range("O30").Value = 6500
$ws.Range("P30").Value = 6250
$ws.Range("R30").Value = 'Región de O''Higgins'
$ws.Range("S30").Value = 3125

# Row 31
$ws.Range("D31").Value = 44454
$ws.Range("M31").Value = 300
$ws.Range("N31").Value = 13000
$ws.Range("O31").Value = 14000
$ws.Range("P31").Value = 13500
$ws.Range("R31").Value = 'Provincia de Limarí'
$ws.Range("S31").Value = 6750
